$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 89.933845520019531
$ws.Range("B3").Value = 87.274131774902344
$ws.Range("B4").Value = 91.155021667480469
$ws.Range("B5").Value = 92.403213500976562
$ws.Range("B6").Value = 100.34432220458984
$ws.Range("B7").Value = 100.99150085449219
$ws.Range("B8").Value = 95.411331176757812
$ws.Range("B9").Value = 99.457733154296875
$ws.Range("B10").Value = 115.80441284179688
$ws.Range("B11").Value = 122.04225921630859
$ws.Range("B12").Value = 106.31224822998047
